# jobInfo.xlsx edit: update row 2 of Sheet1
#  - A2: "java developer" -> "developer"
#  - D2: "0 - 1 an experienta" -> cleared (no value)
#  - selection moves to D5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "developer"
$ws.Range("D2").ClearContents()

$ws.Range("D5").Select()
